$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 16 new rows before row 17 (old Totals row), so data rows 17-32 become available
for ($i = 0; $i -lt 16; $i++) {
    $ws.Rows.Item(17).Insert()
}

# Copy formatting (number formats, fonts, borders, fills, alignment) from row 16 into new rows 17-32
$ws.Range("A16:Q16").Copy()
for ($r = 17; $r -le 32; $r++) {
    $ws.Range("A" + $r + ":Q" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Re-create the merges for each new data row (A:B, C:G, H:K, L:M, N:O)
for ($r = 17; $r -le 32; $r++) {
    $ws.Range("A" + $r + ":B" + $r).Merge()
    $ws.Range("C" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
    $ws.Range("N" + $r + ":O" + $r).Merge()
}

# Set row heights for all data rows + totals + footer to match source layout
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 24.75
$ws.Rows.Item(14).RowHeight = 25.5
$ws.Rows.Item(15).RowHeight = 24.75
$ws.Rows.Item(16).RowHeight = 25.5
$ws.Rows.Item(17).RowHeight = 25.5
$ws.Rows.Item(18).RowHeight = 24.75
$ws.Rows.Item(19).RowHeight = 25.5
$ws.Rows.Item(20).RowHeight = 24.75
$ws.Rows.Item(21).RowHeight = 25.5
$ws.Rows.Item(22).RowHeight = 25.5
$ws.Rows.Item(23).RowHeight = 24.75
$ws.Rows.Item(24).RowHeight = 25.5
$ws.Rows.Item(25).RowHeight = 24.75
$ws.Rows.Item(26).RowHeight = 25.5
$ws.Rows.Item(27).RowHeight = 25.5
$ws.Rows.Item(28).RowHeight = 24.75
$ws.Rows.Item(29).RowHeight = 25.5
$ws.Rows.Item(30).RowHeight = 24.75
$ws.Rows.Item(31).RowHeight = 25.5
$ws.Rows.Item(32).RowHeight = 25.5
$ws.Rows.Item(33).RowHeight = 24.75
$ws.Rows.Item(34).RowHeight = 16.5

# Fill in the data table (26 rows of shortage items)
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = 'ALPHINTERN 30 F.C.TABS'
$ws.Range("H7").Value = '4:0'
$ws.Range("L7").Value = '1'
$ws.Range("N7").Value = '87.00'
$ws.Range("P7").Value = '28.7100'
$ws.Range("Q7").Value = '0:1'

$ws.Range("A8").Value = 2
$ws.Range("C8").Value = 'AVIVAVASC 5/160MG 28 F.C. TAB.'
$ws.Range("H8").Value = '0:0'
$ws.Range("L8").Value = '1'
$ws.Range("N8").Value = '124.00'
$ws.Range("P8").Value = '124.0000'
$ws.Range("Q8").Value = '1:0'

$ws.Range("A9").Value = 3
$ws.Range("C9").Value = 'BI-PROFENID 150MG 20 SCORED TABS.'
$ws.Range("H9").Value = '2:0'
$ws.Range("L9").Value = '1'
$ws.Range("N9").Value = '54.00'
$ws.Range("P9").Value = '27.0000'
$ws.Range("Q9").Value = '0:1'

$ws.Range("A10").Value = 4
$ws.Range("C10").Value = 'CALAMYL LOTION 100 ML'
$ws.Range("H10").Value = '2:0'
$ws.Range("L10").Value = '1'
$ws.Range("N10").Value = '49.00'
$ws.Range("P10").Value = '49.0000'
$ws.Range("Q10").Value = '1:0'

$ws.Range("A11").Value = 5
$ws.Range("C11").Value = 'CETAL 500MG 20 TAB'
$ws.Range("H11").Value = '2:0'
$ws.Range("L11").Value = '1'
$ws.Range("N11").Value = '24.00'
$ws.Range("P11").Value = '12.0000'
$ws.Range("Q11").Value = '0:1'

$ws.Range("A12").Value = 6
$ws.Range("C12").Value = 'CORASORE 150MG/ML ORAL DROPS 15 ML'
$ws.Range("H12").Value = '2:0'
$ws.Range("L12").Value = '1'
$ws.Range("N12").Value = '27.00'
$ws.Range("P12").Value = '27.0000'
$ws.Range("Q12").Value = '1:0'

$ws.Range("A13").Value = 7
$ws.Range("C13").Value = 'DECLOPHEN 75MG/3ML 3 AMPOULES'
$ws.Range("H13").Value = '7:2'
$ws.Range("L13").Value = '1'
$ws.Range("N13").Value = '36.00'
$ws.Range("P13").Value = '11.8800'
$ws.Range("Q13").Value = '0:1'

$ws.Range("A14").Value = 8
$ws.Range("C14").Value = 'FAROVIGA 100MG 12 F.C.TAB.'
$ws.Range("H14").Value = '1:11'
$ws.Range("L14").Value = '1'
$ws.Range("N14").Value = '108.00'
$ws.Range("P14").Value = '35.6400'
$ws.Range("Q14").Value = '0:4'

$ws.Range("A15").Value = 9
$ws.Range("C15").Value = 'FELDENE 20MG/ML I.M. 6 AMP.'
$ws.Range("H15").Value = '2:5'
$ws.Range("L15").Value = '1'
$ws.Range("N15").Value = '63.00'
$ws.Range("P15").Value = '63.0000'
$ws.Range("Q15").Value = '1:0'

$ws.Range("A16").Value = 10
$ws.Range("C16").Value = 'FUSI 2% CREAM 15 GM'
$ws.Range("H16").Value = '2:0'
$ws.Range("L16").Value = '1'
$ws.Range("N16").Value = '35.00'
$ws.Range("P16").Value = '35.0000'
$ws.Range("Q16").Value = '1:0'

$ws.Range("A17").Value = 11
$ws.Range("C17").Value = 'LOADLESS 5/20MG 30 CAP'
$ws.Range("H17").Value = '1:0'
$ws.Range("L17").Value = '1'
$ws.Range("N17").Value = '114.00'
$ws.Range("P17").Value = '37.6200'
$ws.Range("Q17").Value = '0:1'

$ws.Range("A18").Value = 12
$ws.Range("C18").Value = 'MAXOPHAGE XR 1000MG 30 EXT. REL. TABS.'
$ws.Range("H18").Value = '2:1'
$ws.Range("L18").Value = '1'
$ws.Range("N18").Value = '72.00'
$ws.Range("P18").Value = '23.7600'
$ws.Range("Q18").Value = '0:1'

$ws.Range("A19").Value = 13
$ws.Range("C19").Value = 'MELANOFREE CREAM 30 GM'
$ws.Range("H19").Value = '2:0'
$ws.Range("L19").Value = '1'
$ws.Range("N19").Value = '74.00'
$ws.Range("P19").Value = '74.0000'
$ws.Range("Q19").Value = '1:0'

$ws.Range("A20").Value = 14
$ws.Range("C20").Value = 'MUCO 15MG/5ML SYRUP 100ML'
$ws.Range("H20").Value = '1:0'
$ws.Range("L20").Value = '1'
$ws.Range("N20").Value = '35.00'
$ws.Range("P20").Value = '35.0000'
$ws.Range("Q20").Value = '1:0'

$ws.Range("A21").Value = 15
$ws.Range("C21").Value = 'NASSAR 12*8 F.C. TAB.'
$ws.Range("H21").Value = '1:10'
$ws.Range("L21").Value = '1'
$ws.Range("N21").Value = '156.00'
$ws.Range("P21").Value = '12.4800'
$ws.Range("Q21").Value = '0:1'

$ws.Range("A22").Value = 16
$ws.Range("C22").Value = 'OCTOVENT PLUS SYRUP 100 ML'
$ws.Range("H22").Value = '2:0'
$ws.Range("L22").Value = '1'
$ws.Range("N22").Value = '29.00'
$ws.Range("P22").Value = '29.0000'
$ws.Range("Q22").Value = '1:0'

$ws.Range("A23").Value = 17
$ws.Range("C23").Value = 'ORS 10 SACHET'
$ws.Range("H23").Value = '6:9'
$ws.Range("L23").Value = '1'
$ws.Range("N23").Value = '40.00'
$ws.Range("P23").Value = '4.0000'
$ws.Range("Q23").Value = '0:1'

$ws.Range("A24").Value = 18
$ws.Range("C24").Value = 'PHYTO K 10 MG 50 F.C.TAB.'
$ws.Range("H24").Value = '1:0'
$ws.Range("L24").Value = '1'
$ws.Range("N24").Value = '72.50'
$ws.Range("P24").Value = '-14.5000'
$ws.Range("Q24").Value = '0:-1'

$ws.Range("A25").Value = 19
$ws.Range("C25").Value = 'PREDSOL 5 MG /5ML  SYRUP'
$ws.Range("H25").Value = '1:0'
$ws.Range("L25").Value = '0'
$ws.Range("N25").Value = '57.00'
$ws.Range("P25").Value = '57.0000'
$ws.Range("Q25").Value = '1:0'

$ws.Range("A26").Value = 20
$ws.Range("C26").Value = 'PULMICORT 0.25MG/ML 20 NEBULIZER VIAL SUSP.'
$ws.Range("H26").Value = '0:9'
$ws.Range("L26").Value = '1'
$ws.Range("N26").Value = '564.00'
$ws.Range("P26").Value = '28.2000'
$ws.Range("Q26").Value = '0:1'

$ws.Range("A27").Value = 21
$ws.Range("C27").Value = 'RIVO 320MG 20*10 TABS'
$ws.Range("H27").Value = '0:11'
$ws.Range("L27").Value = '1'
$ws.Range("N27").Value = '141.00'
$ws.Range("P27").Value = '7.0500'
$ws.Range("Q27").Value = '0:1'

$ws.Range("A28").Value = 22
$ws.Range("C28").Value = 'UNICTAM 750 MG I.M/I.V VIAL'
$ws.Range("H28").Value = '6:0'
$ws.Range("L28").Value = '1'
$ws.Range("N28").Value = '39.00'
$ws.Range("P28").Value = '39.0000'
$ws.Range("Q28").Value = '1:0'

$ws.Range("A29").Value = 23
$ws.Range("C29").Value = 'VOLTAREN 75MG/3ML 3 AMP.'
$ws.Range("H29").Value = '1:0'
$ws.Range("L29").Value = '1'
$ws.Range("N29").Value = '51.00'
$ws.Range("P29").Value = '33.6600'
$ws.Range("Q29").Value = '0:2'

$ws.Range("A30").Value = 24
$ws.Range("C30").Value = 'سرنجات 3 سم'
$ws.Range("H30").Value = '0:0'
$ws.Range("L30").Value = '0'
$ws.Range("N30").Value = '2.00'
$ws.Range("P30").Value = '14.0000'
$ws.Range("Q30").Value = '7:0'

$ws.Range("A31").Value = 25
$ws.Range("C31").Value = 'سرنجات 5 سم'
$ws.Range("H31").Value = '0:0'
$ws.Range("L31").Value = '0'
$ws.Range("N31").Value = '3.00'
$ws.Range("P31").Value = '21.0000'
$ws.Range("Q31").Value = '7:0'

$ws.Range("A32").Value = 26
$ws.Range("C32").Value = 'كالونا '
$ws.Range("H32").Value = '0:0'
$ws.Range("L32").Value = '0'
$ws.Range("N32").Value = '15.00'
$ws.Range("P32").Value = '15.0000'
$ws.Range("Q32").Value = '1:0'

# Totals row (now row 33)
$ws.Range("P33").Value = 829.5

# Footer row (now row 34): update generation timestamp
$ws.Range("A34").Value = 'Sunday, 15 June, 2025 1:41 PM'
